$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B16").Value = "Implement a notification to tell the user that edits have been successfully saved."
$ws.Range("C16").Value = "Done"

$ws.Range("B17").Select()
